# Update "想去人数" (interest count, column F) for a set of events whose
# data is duplicated across the "展览" (sheet1) and "全部类型" (sheet4) sheets.
# The mapping below goes from worksheet name -> { row number -> new value }.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        3  = 1212
        4  = 14583
        5  = 17592
        8  = 55
        16 = 39
        17 = 143
        19 = 1329
        23 = 212
        24 = 7233
        28 = 1172
        30 = 5860
        33 = 137
        35 = 223
        36 = 5079
    }
    "全部类型" = @{
        3  = 1212
        4  = 14583
        5  = 17592
        8  = 55
        16 = 39
        17 = 143
        19 = 1329
        24 = 212
        25 = 7233
        29 = 1172
        32 = 5860
        35 = 137
        37 = 223
        38 = 5079
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $updates[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowMap[$row]
    }
}
